# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new column headers, matching the style of the
# existing header cells (bold, centered header style = style index used
# by the other header cells, exposed via Style "Header" is not named, so
# copy formatting from an existing header cell instead).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the rest of the header row (bold + centered + bordered).
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A1").Select() | Out-Null

# Season record values for every data row (rows 2 through 45).
$wins = 74
$losses = 88
$ties = 0

for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
